{"js": "// The SDO template currently builds the case-management venue line from\n// `caseManagementLocation.venue_name`. Switch it to use the location ref\n// data's external short name field instead: `caseManagementLocation.external_short_name`.\nconst body = context.document.body;\n\nconst results = body.search(\"venue_name\", { matchCase: true, matchWholeWord: false });\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find 'venue_name' placeholder text to update.\");\n}\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(\"external_short_name\", Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# The SDO template currently builds the case-management venue line from\n# `caseManagementLocation.venue_name`. Switch it to use the location ref\n# data's external short name field instead:\n# `caseManagementLocation.external_short_name`.\n\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n\n$find.Execute(\n    \"venue_name\",          # FindText\n    $true,                  # MatchCase\n    $false,                 # MatchWholeWord\n    $false,                 # MatchWildcards\n    $false,                 # MatchSoundsLike\n    $false,                 # MatchAllWordForms\n    $true,                  # Forward\n    1,                      # Wrap (wdFindContinue)\n    $false,                 # Format\n    \"external_short_name\",  # ReplaceWith\n    2                       # Replace (wdReplaceAll)\n) | Out-Null\n"}
